$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.109299999999996
$ws.Range("A9").Value = -20.39239999999998
$ws.Range("A18").Value = -22.95460000000002
$ws.Range("A20").Value = -22.09170000000002
